# parental-relationship-types.xlsx update
# New version of the controlled vocabulary: header renamed, a "note" column
# added to flag deprecated ("Non Valido") items, a couple of label
# capitalization fixes, and the "Intestatario Scheda" label expanded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "codice_1_livello"
$ws.Range("B1").Value = "label_1_livello_it"
$ws.Range("C1").Value = "note"

# Header used to be bold (style index 1 on A1/B1); now it is plain text.
$ws.Range("A1:C1").Font.Bold = $false

# --- Mark deprecated rows with "Non Valido" in the new note column --------
$deprecatedRows = @(15, 16, 20, 27, 29, 30, 31)
foreach ($r in $deprecatedRows) {
    $ws.Cells.Item($r, 3).Value = "Non Valido"
}

# --- Row 2: expand "Intestatario Scheda" label, drop its special style ----
$ws.Range("B2").Value = "Intestatario della Scheda"
$ws.Range("B2").Font.Bold = $false

# --- Capitalization fixes ---------------------------------------------------
$ws.Range("B12").Value = "Zio / Zia (collaterale)"
$ws.Range("B21").Value = "Nipote (affine)"
$ws.Range("B22").Value = "Zio / Zia (affine)"

# --- View bookkeeping (matches where the author ended up in the sheet) ----
$ws.Range("C32").Select()
$ws.Application.ActiveWindow.ScrollRow = 20
